$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column before the "Bemerkung" column (L) for the new "Ziel" column ---
$ws.Columns("L").Insert()
$ws.Range("L1").Value = "Ziel"
$ws.Range("L2").Value = 5
$ws.Range("L3").Value = 2

# --- Insert a brand-new tracked poker session as row 4 (shifts old rows down) ---
$ws.Rows(4).Insert()
$ws.Range("A3:N3").Copy()
$ws.Range("A4:N4").PasteSpecial(-4122)

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "sc.ch"
$ws.Range("C4").Value = 0.8
$ws.Range("D4").Value = 0.01
$ws.Range("E4").Value = 45920
$ws.Range("F4").Value = 0.91666666666666663
$ws.Range("G4").Value = 45920
$ws.Range("H4").Value = 0.97916666666666663
$ws.Range("I4").Value = 1.2
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1.45
$ws.Range("L4").Value = 2
$ws.Range("M4").Value = "AA Slowplay war Fehler, Flush des Gegners nicht wahrgenommen"
$ws.Range("N4").Value = "Continuationbets waren Killer, generell nichts getroffen auch bei Flush/Straightdraws, Es muss bereits bei Beginn weniger gespielt werden."

# --- Append a new trailing blank tracked row (row 42) ---
$ws.Range("A41:N41").Copy()
$ws.Range("A42:N42").PasteSpecial(-4122)
$ws.Range("A42").Value = 40

# --- Restore the previously-selected cell ---
$ws.Range("J20").Select()
